$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Rectangle 65 (shape id 4 / shapes index 1) - widen/reposition the outer
# sequence-diagram frame (moved left edge further left, increased width).
$rect65 = $s.Shapes.Item(1)
$rect65.Left = 44.399845123291016
$rect65.Top = 37.80299377441406
$rect65.Width = 787.3932495117188
$rect65.Height = 471.52276611328125

# Straight Arrow Connector 15 (shape id 16 / shapes index 14) - matching
# lifeline connector adjusted to follow the resized frame above.
$conn15 = $s.Shapes.Item(14)
$conn15.Left = 41.733150482177734
$conn15.Top = 480.6029968261719
$conn15.Width = 212.04623413085938
$conn15.Height = 0.8138582706451416

# Rectangle 28 (shape id 29 / shapes index 25) - merge the split "r" /
# "esult :Command" runs back into a single "result :Command" run.
$rect28 = $s.Shapes.Item(25)
$rect28.TextFrame.TextRange.Characters(1, 15).Text = "result :Command"
